# Updating filtered feeds from workflow
# Append two new feed rows (69 and 70) to the "Filtered Feeds" worksheet,
# matching the existing link/keywords/title pattern used by the rest of
# the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$link69 = "https://www.genomeweb.com/cancer/fda-proposal-reclassify-cdx-assays-may-broaden-opportunities-dx-manufacturers-experts-say"
$link70 = "https://www.360dx.com/cancer/fda-proposal-reclassify-cdx-assays-may-broaden-opportunities-dx-manufacturers-experts-say"
$keywords = "CDx, companion diagnostic"
$title = "FDA Proposal to Reclassify CDx Assays May Broaden Opportunities for Dx Manufacturers, Experts Say"

# Add the hyperlinks first (creates the relationship + <hyperlinks> entries),
# then reapply the same hyperlink-cell style used throughout column A so we
# don't leave the cells on a freshly auto-generated style.
$ws.Hyperlinks.Add($ws.Range("A69"), $link69)
$ws.Hyperlinks.Add($ws.Range("A70"), $link70)
$ws.Range("A69").Style = $ws.Range("A2").Style
$ws.Range("A70").Style = $ws.Range("A2").Style

$ws.Range("B69").Value = $keywords
$ws.Range("C69").Value = $title

$ws.Range("B70").Value = $keywords
$ws.Range("C70").Value = $title
